# Generate Report for Handoff
# Adds two new rows (3 and 4) to the Overview / zh-cn / de-de sheets,
# describing two additional files that went through handoff, and
# refreshes the handoff timestamp / file name for the existing row.
#
# NOTE: this engine's `Range(...).Hyperlinks.Delete()` clears the *entire*
# worksheet hyperlink collection rather than just the hyperlinks that
# intersect the given range, so every sheet's hyperlinks are dropped once
# up front and then fully re-added (old + new) in final left-to-right,
# top-to-bottom order.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Source data
# ---------------------------------------------------------------------

$mdRepo  = "https://github.com/OpenLocalizationTest/oltest/blob/a97e9ff0a7f65c484d2851797170f8f1da6b4474/e2e/"
$zhRepo  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0fbe76e7d0cfdec267c4716cd38b0f437b4dd0db/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$deRepo  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6131e1dae92dad3fc00f37e247e113a5293a1b62/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

# Source file names for the three handed-off files (row2 = existing / refreshed, row3, row4)
$srcName2 = "1cdaebeb-e694-402f-9a9d-8ab567d74374.png"
$srcName3 = "c12829f7-f2d0-4a6c-b335-59d99ebc9d46.md"
$srcName4 = "ea797b38-0867-418d-bf86-0fcb5489d0c7.png"

$ext2 = ".png"
$ext3 = ".md"
$ext4 = ".png"

$target2 = "3d1ae97d79a99e10600e226a5f7ab0c317e5d818.png"
$target3zh = "c12829f7-f2d0-4a6c-b335-59d99ebc9d46.350ef275c1165c7d758dfdd1c2bcf9f84b1e90d7.zh-cn.xlf"
$target3de = "c12829f7-f2d0-4a6c-b335-59d99ebc9d46.350ef275c1165c7d758dfdd1c2bcf9f84b1e90d7.de-de.xlf"
$target4 = "70521efce71e9926578fa7d9a8e4820974750e16.png"

$status = "Ready for handoff"
$overviewDate = "2016-33-18 04:33:08"
$zhDatetime = "2016-03-18 04:33:05"
$deDatetime = "2016-03-18 04:33:08"
$epoch = "0001-01-01 00:00:00"

$dependencyFrom = "e2e\c12829f7-f2d0-4a6c-b335-59d99ebc9d46.md"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================

$ws1 = $wb.Worksheets.Item("Overview")

# --- refresh row 2 (existing file renamed from .md to .png) ---
$ws1.Range("D2").Value2 = $overviewDate

# --- row 3 (new file) ---
$ws1.Cells.Item(3, 2).Value2 = $status
$ws1.Cells.Item(3, 3).Value2 = $status
$ws1.Cells.Item(3, 4).Value2 = $overviewDate

# --- row 4 (new file) ---
$ws1.Cells.Item(4, 2).Value2 = $status
$ws1.Cells.Item(4, 3).Value2 = $status
$ws1.Cells.Item(4, 4).Value2 = $overviewDate

# drop all existing hyperlinks, then re-add every hyperlink (old + new) in order
$ws1.Hyperlinks.Delete()
$null = $ws1.Hyperlinks.Add($ws1.Range("A2"), ($mdRepo + $srcName2), "", "", $srcName2)
$null = $ws1.Hyperlinks.Add($ws1.Range("A3"), ($mdRepo + $srcName3), "", "", $srcName3)
$null = $ws1.Hyperlinks.Add($ws1.Range("A4"), ($mdRepo + $srcName4), "", "", $srcName4)

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================

$ws2 = $wb.Worksheets.Item("zh-cn")

# --- refresh row 2 ---
$ws2.Cells.Item(2, 5).Value2 = $zhDatetime
$ws2.Cells.Item(2, 8).Value2 = $epoch
$ws2.Cells.Item(2, 9).Value2 = "IsDependency"
$ws2.Cells.Item(2, 10).Value2 = $dependencyFrom

# --- row 3 (new file, Include reason, no dependency) ---
$ws2.Cells.Item(3, 3).Value2 = $status
$ws2.Cells.Item(3, 5).Value2 = $zhDatetime
$ws2.Cells.Item(3, 8).Value2 = $epoch
$ws2.Cells.Item(3, 9).Value2 = "Include"

# --- row 4 ---
$ws2.Cells.Item(4, 3).Value2 = $status
$ws2.Cells.Item(4, 5).Value2 = $zhDatetime
$ws2.Cells.Item(4, 8).Value2 = $epoch
$ws2.Cells.Item(4, 9).Value2 = "IsDependency"
$ws2.Cells.Item(4, 10).Value2 = $dependencyFrom

# apply the same date/time style used in E2 to the new E3/E4 cells
$ws2.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# drop all existing hyperlinks, then re-add every hyperlink (old + new) in order
$ws2.Hyperlinks.Delete()
$null = $ws2.Hyperlinks.Add($ws2.Range("A2"), ($mdRepo + $srcName2), "", "", $srcName2)
$null = $ws2.Hyperlinks.Add($ws2.Range("B2"), ($mdRepo + $srcName2), "", "", $ext2)
$null = $ws2.Hyperlinks.Add($ws2.Range("D2"), ($zhRepo + $target2), "", "", $target2)
$null = $ws2.Hyperlinks.Add($ws2.Range("A3"), ($mdRepo + $srcName3), "", "", $srcName3)
$null = $ws2.Hyperlinks.Add($ws2.Range("B3"), ($mdRepo + $srcName3), "", "", $ext3)
$null = $ws2.Hyperlinks.Add($ws2.Range("D3"), ($zhRepo + $target3zh), "", "", $target3zh)
$null = $ws2.Hyperlinks.Add($ws2.Range("A4"), ($mdRepo + $srcName4), "", "", $srcName4)
$null = $ws2.Hyperlinks.Add($ws2.Range("B4"), ($mdRepo + $srcName4), "", "", $ext4)
$null = $ws2.Hyperlinks.Add($ws2.Range("D4"), ($zhRepo + $target4), "", "", $target4)

# ===========================================================================
# Sheet "de-de"
# ===========================================================================

$ws3 = $wb.Worksheets.Item("de-de")

# --- refresh row 2 ---
$ws3.Cells.Item(2, 5).Value2 = $deDatetime
$ws3.Cells.Item(2, 8).Value2 = $epoch
$ws3.Cells.Item(2, 9).Value2 = "IsDependency"
$ws3.Cells.Item(2, 10).Value2 = $dependencyFrom

# --- row 3 ---
$ws3.Cells.Item(3, 3).Value2 = $status
$ws3.Cells.Item(3, 5).Value2 = $deDatetime
$ws3.Cells.Item(3, 8).Value2 = $epoch
$ws3.Cells.Item(3, 9).Value2 = "Include"

# --- row 4 ---
$ws3.Cells.Item(4, 3).Value2 = $status
$ws3.Cells.Item(4, 5).Value2 = $deDatetime
$ws3.Cells.Item(4, 8).Value2 = $epoch
$ws3.Cells.Item(4, 9).Value2 = "IsDependency"
$ws3.Cells.Item(4, 10).Value2 = $dependencyFrom

# apply the same date/time style used in E2 to the new E3/E4 cells
$ws3.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# drop all existing hyperlinks, then re-add every hyperlink (old + new) in order
$ws3.Hyperlinks.Delete()
$null = $ws3.Hyperlinks.Add($ws3.Range("A2"), ($mdRepo + $srcName2), "", "", $srcName2)
$null = $ws3.Hyperlinks.Add($ws3.Range("B2"), ($mdRepo + $srcName2), "", "", $ext2)
$null = $ws3.Hyperlinks.Add($ws3.Range("D2"), ($deRepo + $target2), "", "", $target2)
$null = $ws3.Hyperlinks.Add($ws3.Range("A3"), ($mdRepo + $srcName3), "", "", $srcName3)
$null = $ws3.Hyperlinks.Add($ws3.Range("B3"), ($mdRepo + $srcName3), "", "", $ext3)
$null = $ws3.Hyperlinks.Add($ws3.Range("D3"), ($deRepo + $target3de), "", "", $target3de)
$null = $ws3.Hyperlinks.Add($ws3.Range("A4"), ($mdRepo + $srcName4), "", "", $srcName4)
$null = $ws3.Hyperlinks.Add($ws3.Range("B4"), ($mdRepo + $srcName4), "", "", $ext4)
$null = $ws3.Hyperlinks.Add($ws3.Range("D4"), ($deRepo + $target4), "", "", $target4)

Write-Host "Handoff report rows added."
